$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Comparison")

# Updated auto-increase projection values (new-hire-only and all-eligible scenarios)

# Row 2
$ws.Range("E2").Value = 82
$ws.Range("F2").Value = 0.803921568627451
$ws.Range("G2").Value = 0.803921568627451
$ws.Range("H2").Value = 0.1020570907636317
$ws.Range("I2").Value = 0.08204589649625291
$ws.Range("J2").Value = 452494.6720494052
$ws.Range("K2").Value = 163633.3362257026
$ws.Range("M2").Value = 163633.3362257026
$ws.Range("N2").Value = 616128.0082751078
$ws.Range("O2").Value = 10292112.5688
$ws.Range("P2").Value = 9884371.638699999
$ws.Range("Q2").Value = 0.01589890657839757
$ws.Range("R2").Value = 0.01655475352474948

# Row 3
$ws.Range("D3").Value = 102
$ws.Range("E3").Value = 87
$ws.Range("F3").Value = 0.8529411764705882
$ws.Range("G3").Value = 0.8446601941747572
$ws.Range("H3").Value = 0.09750327982997359
$ws.Range("I3").Value = 0.08235713927386119
$ws.Range("J3").Value = 473488.5336105639
$ws.Range("K3").Value = 171531.8470123119
$ws.Range("M3").Value = 171531.8470123119
$ws.Range("N3").Value = 645020.3806228759
$ws.Range("O3").Value = 10436953.712764
$ws.Range("P3").Value = 10029580.554761
$ws.Range("Q3").Value = 0.01643504912765254
$ws.Range("R3").Value = 0.01710259427856995

# Row 4
$ws.Range("E4").Value = 87
$ws.Range("F4").Value = 0.8365384615384616
$ws.Range("G4").Value = 0.8365384615384616
$ws.Range("H4").Value = 0.09749531305237304
$ws.Range("I4").Value = 0.08155857918804287
$ws.Range("J4").Value = 497757.4148347
$ws.Range("K4").Value = 176579.9150305909
$ws.Range("M4").Value = 176579.9150305909
$ws.Range("N4").Value = 674337.3298652909
$ws.Range("O4").Value = 10911311.62264692
$ws.Range("P4").Value = 10502967.26990383
$ws.Range("Q4").Value = 0.01618319787183892
$ws.Range("R4").Value = 0.01681238363339275

# Row 5
$ws.Range("E5").Value = 86
$ws.Range("F5").Value = 0.819047619047619
$ws.Range("G5").Value = 0.819047619047619
$ws.Range("H5").Value = 0.09816608625862366
$ws.Range("I5").Value = 0.08040269922134892
$ws.Range("J5").Value = 513130.3226945847
$ws.Range("K5").Value = 181692.4051889305
$ws.Range("M5").Value = 181692.4051889305
$ws.Range("N5").Value = 694822.7278835152
$ws.Range("O5").Value = 11330655.55062633
$ws.Range("P5").Value = 10919960.86730095
$ws.Range("Q5").Value = 0.01603547159095195
$ws.Range("R5").Value = 0.01663855826928791

# Row 6
$ws.Range("E6").Value = 88
$ws.Range("F6").Value = 0.8301886792452831
$ws.Range("G6").Value = 0.8301886792452831
$ws.Range("H6").Value = 0.09702974602399582
$ws.Range("I6").Value = 0.08055299669916634
$ws.Range("J6").Value = 535915.2845663553
$ws.Range("K6").Value = 190023.7034400649
$ws.Range("M6").Value = 190023.7034400649
$ws.Range("N6").Value = 725938.9880064201
$ws.Range("O6").Value = 11688284.31584512
$ws.Range("P6").Value = 11273818.79201997
$ws.Range("Q6").Value = 0.01625762158971963
$ws.Range("R6").Value = 0.01685530936283726

# Row 7
$ws.Range("D7").Value = 102
$ws.Range("E7").Value = 88
$ws.Range("F7").Value = 0.8627450980392157
$ws.Range("G7").Value = 0.8627450980392157
$ws.Range("H7").Value = 0.09833450573160141
$ws.Range("I7").Value = 0.08483761278804831
$ws.Range("J7").Value = 472190.6408301356
$ws.Range("K7").Value = 173481.3206160678
$ws.Range("M7").Value = 173481.3206160678
$ws.Range("N7").Value = 645671.9614462035
$ws.Range("O7").Value = 10262129.6088
$ws.Range("P7").Value = 9854388.678699998
$ws.Range("Q7").Value = 0.01690500190791819
$ws.Range("R7").Value = 0.01760447312079775

# Row 8
$ws.Range("D8").Value = 103
$ws.Range("E8").Value = 89
$ws.Range("F8").Value = 0.8640776699029126
$ws.Range("G8").Value = 0.8640776699029126
$ws.Range("H8").Value = 0.09614705683641531
$ws.Range("I8").Value = 0.08307852483923264
$ws.Range("J8").Value = 477788.110972513
$ws.Range("K8").Value = 173681.6356932865
$ws.Range("M8").Value = 173681.6356932865
$ws.Range("N8").Value = 651469.7466657992
$ws.Range("O8").Value = 10483190.267664
$ws.Range("P8").Value = 10075817.109661
$ws.Range("Q8").Value = 0.01656763172838877
$ws.Range("R8").Value = 0.01723747402349683

# Row 9
$ws.Range("D9").Value = 103
$ws.Range("E9").Value = 89
$ws.Range("F9").Value = 0.8640776699029126
$ws.Range("G9").Value = 0.8557692307692307
$ws.Range("H9").Value = 0.0970994315717198
$ws.Range("I9").Value = 0.08319085971041408
$ws.Range("J9").Value = 509184.1566766572
$ws.Range("K9").Value = 182293.2859515694
$ws.Range("M9").Value = 182293.2859515694
$ws.Range("N9").Value = 691477.4426282267
$ws.Range("O9").Value = 10842468.31879392
$ws.Range("P9").Value = 10434123.96605083
$ws.Range("Q9").Value = 0.01681289542120121
$ws.Range("R9").Value = 0.01747087599732293

# Row 10
$ws.Range("E10").Value = 91
$ws.Range("F10").Value = 0.8666666666666667
$ws.Range("G10").Value = 0.8666666666666667
$ws.Range("H10").Value = 0.09511055620966978
$ws.Range("I10").Value = 0.08242914871504717
$ws.Range("J10").Value = 526523.2199554271
$ws.Range("K10").Value = 188388.8538193516
$ws.Range("M10").Value = 188388.8538193516
$ws.Range("N10").Value = 714912.0737747787
$ws.Range("O10").Value = 11280831.86545774
$ws.Range("P10").Value = 10870137.18213236
$ws.Range("Q10").Value = 0.01669990795591983
$ws.Range("R10").Value = 0.01733086258828576

# Row 11
$ws.Range("E11").Value = 92
$ws.Range("F11").Value = 0.8679245283018868
$ws.Range("G11").Value = 0.8679245283018868
$ws.Range("H11").Value = 0.09510156385804072
$ws.Range("I11").Value = 0.08254097995226177
$ws.Range("J11").Value = 551974.7101267558
$ws.Range("K11").Value = 198053.4162202651
$ws.Range("M11").Value = 198053.4162202651
$ws.Range("N11").Value = 750028.1263470209
$ws.Range("O11").Value = 11729815.62832147
$ws.Range("P11").Value = 11315350.10449633
$ws.Range("Q11").Value = 0.01688461460059679
$ws.Range("R11").Value = 0.01750307453072668

# Row 12
$ws.Range("H12").Value = 0.1065069431212171
$ws.Range("I12").Value = 0.09084415736809691
$ws.Range("J12").Value = 528276.0912988938
$ws.Range("K12").Value = 201524.0458504469
$ws.Range("M12").Value = 201524.0458504469
$ws.Range("N12").Value = 729800.1371493406
$ws.Range("O12").Value = 10161399.9688
$ws.Range("P12").Value = 9753659.038699998
$ws.Range("Q12").Value = 0.01983231114504054
$ws.Range("R12").Value = 0.02066137898104204

# Row 13
$ws.Range("D13").Value = 103
$ws.Range("E13").Value = 88
$ws.Range("F13").Value = 0.8543689320388349
$ws.Range("G13").Value = 0.8543689320388349
$ws.Range("H13").Value = 0.1101170774262103
$ws.Range("I13").Value = 0.094080609839869
$ws.Range("J13").Value = 595699.0697549444
$ws.Range("K13").Value = 232637.1150845023
$ws.Range("M13").Value = 232637.1150845023
$ws.Range("N13").Value = 828336.1848394468
$ws.Range("O13").Value = 10536239.823564
$ws.Range("P13").Value = 10128866.665561
$ws.Range("Q13").Value = 0.02207970955294848
$ws.Range("R13").Value = 0.0229677339791122

# Row 14
$ws.Range("D14").Value = 103
$ws.Range("E14").Value = 87
$ws.Range("F14").Value = 0.8446601941747572
$ws.Range("G14").Value = 0.8365384615384616
$ws.Range("H14").Value = 0.1155345313075407
$ws.Range("I14").Value = 0.09664907907457733
$ws.Range("J14").Value = 659315.5951792673
$ws.Range("K14").Value = 257359.0052028746
$ws.Range("M14").Value = 257359.0052028746
$ws.Range("N14").Value = 916674.600382142
$ws.Range("O14").Value = 10860875.61657092
$ws.Range("P14").Value = 10452531.26382783
$ws.Range("Q14").Value = 0.02369597206418703
$ws.Range("R14").Value = 0.02462169198129975

# Row 15
$ws.Range("E15").Value = 85
$ws.Range("F15").Value = 0.8095238095238095
$ws.Range("G15").Value = 0.8095238095238095
$ws.Range("H15").Value = 0.1172245253426888
$ws.Range("I15").Value = 0.09489604432503375
$ws.Range("J15").Value = 678472.9241290929
$ws.Range("K15").Value = 264363.7059061846
$ws.Range("M15").Value = 264363.7059061846
$ws.Range("N15").Value = 942836.6300352775
$ws.Range("O15").Value = 11360984.84936805
$ws.Range("P15").Value = 10950290.16604267
$ws.Range("Q15").Value = 0.02326943565292139
$ws.Range("R15").Value = 0.02414216444473664

# Row 16
$ws.Range("D16").Value = 105
$ws.Range("F16").Value = 0.8571428571428571
$ws.Range("H16").Value = 0.1134076443070747
$ws.Range("I16").Value = 0.09628950931732752
$ws.Range("J16").Value = 716159.0169867697
$ws.Range("K16").Value = 280145.5696502721
$ws.Range("M16").Value = 280145.5696502721
$ws.Range("N16").Value = 996304.5866370418
$ws.Range("O16").Value = 11698901.34244909
$ws.Range("P16").Value = 11284435.81862395
$ws.Range("Q16").Value = 0.02394631439738472
$ws.Range("R16").Value = 0.02482583747677638
